# Regenerate the "K" column (column G) values for the save_data sheet.
# The commit regenerates this save-data workbook so the K column (header "K",
# formerly computed from a "Strike#" source) is recalculated/rewritten with
# new values for rows 2-80 (row 1 is the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of worksheet row number -> new value for column G ("K").
# Only rows whose K value actually changes are listed; rows not present here
# keep their existing value (they were already correct / unaffected).
$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 2
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 2
    26 = 2
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 2
    33 = 2
    34 = 1
    35 = 2
    36 = 0
    37 = 1
    38 = 0
    39 = 2
    40 = 0
    41 = 1
    42 = 0
    43 = 0
    44 = 3
    45 = 1
    46 = 1
    47 = 1
    48 = 2
    49 = 0
    50 = 1
    51 = 2
    52 = 2
    54 = 2
    55 = 1
    56 = 2
    57 = 1
    58 = 3
    60 = 1
    62 = 1
    63 = 1
    65 = 2
    66 = 2
    67 = 2
    68 = 1
    69 = 1
    70 = 1
    71 = 0
    72 = 2
    73 = 3
    74 = 1
    75 = 2
    76 = 1
    77 = 2
    79 = 1
    80 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
